$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 15.8388614654541
$ws.Range("D2").Value = 123

$ws.Range("C3").Value = 15.21706581115723
$ws.Range("D3").Value = 123

$ws.Range("C4").Value = 15.96903800964355
$ws.Range("D4").Value = 123

$ws.Range("C5").Value = 16.1430835723877
$ws.Range("D5").Value = 123

$ws.Range("C6").Value = 16.71004295349121
$ws.Range("D6").Value = 123
